# actualizar path reproducible
# Update the tale_path column in the "datos_audio_reproducible" sheet:
# rename the folder segment "AudiosReproducible" -> "AudiosReproducibles"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("datos_audio_reproducible")

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val.Contains("RepositorioDatosProyectoAP/AudiosReproducible/")) {
        $cell.Value2 = $val.Replace("RepositorioDatosProyectoAP/AudiosReproducible/", "RepositorioDatosProyectoAP/AudiosReproducibles/")
    }
}
